$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 234, shifting existing rows 234-279 down to 235-280.
# (Excel copies the formatting of the row above into the newly inserted row,
# matching the original cell style "s=2" used throughout this column.)
$ws.Rows.Item(234).Insert()

# Populate the new row with the new category entry
$ws.Cells.Item(234, 1).Value = "on-09-07"
$ws.Cells.Item(234, 2).Value = "Trop-2"

# Update sheet view to match target state
$ws.Application.ActiveWindow.ScrollRow = 300
$ws.Range("B224").Select() | Out-Null
